$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.669.05'
$ws.Range('E2').Value = '  -0.65%  '
$ws.Range('D3').Value = '3.842.62'
$ws.Range('E3').Value = '  +2.77%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''601.21'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').Value = '''163.84'
$ws.Range('E6').Value = '  -2.08%  '
$ws.Range('D7').Value = '3.841.04'
$ws.Range('E7').Value = '  +2.80%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('E9').Value = '  -2.03%  '
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('D11').Value = '''6.35'
$ws.Range('E11').Value = '  -0.34%  '
$ws.Range('E12').Value = '  -0.09%  '
$ws.Range('D13').Value = '''36.96'
$ws.Range('E13').Value = '  -3.09%  '
$ws.Range('D14').Value = '''0.0000244'
$ws.Range('E14').Value = '  -1.14%  '
$ws.Range('D15').Value = '4.489.90'
$ws.Range('E15').Value = '  +2.90%  '
$ws.Range('D16').Value = '3.863.09'
$ws.Range('E16').Value = '  +3.28%  '
$ws.Range('D17').Value = '68.819.33'
$ws.Range('E17').Value = '  -0.32%  '
$ws.Range('D18').Value = '''7.56'
$ws.Range('E18').Value = '  +2.68%  '
$ws.Range('E19').Value = '  -0.35%  '
$ws.Range('D20').Value = '''17.15'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').Value = '''11.22'
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '''486.39'
$ws.Range('E22').Value = '  -1.10%  '
$ws.Range('E23').Value = '  -1.28%  '
$ws.Range('D24').Value = '''0.0000159'
$ws.Range('E24').Value = '  +6.04%  '
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('D26').Value = '''2.24'
$ws.Range('E26').Value = '  -2.04%  '
$ws.Range('D27').Value = '''12.12'
$ws.Range('E27').Value = '  -1.27%  '
$ws.Range('E28').Value = '  -0.44%  '
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.53%  '
$ws.Range('D31').Value = '''7.86'
$ws.Range('E31').Value = '  -3.95%  '
$ws.Range('D32').Value = '3.996.30'
$ws.Range('E32').Value = '  +2.91%  '
$ws.Range('D33').Value = '''2.37'
$ws.Range('E33').Value = '  -3.66%  '
$ws.Range('D34').Value = '''31.96'
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('D35').Value = '3.790.41'
$ws.Range('E35').Value = '  +3.25%  '
$ws.Range('E36').Value = '  -1.12%  '
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('E38').Value = '  +0.72%  '
$ws.Range('E39').Value = '  -0.75%  '
$ws.Range('D40').Value = '''1.00'
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('D42').Value = '''2.98'
$ws.Range('E42').Value = '  -1.84%  '
$ws.Range('D43').Value = '''432.23'
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  -0.65%  '
$ws.Range('E45').Value = '  +0.25%  '
$ws.Range('D47').Value = '''8.42'
$ws.Range('E47').Value = '  -0.59%  '
$ws.Range('D48').Value = '2.843.53'
$ws.Range('E48').Value = '  +2.30%  '
$ws.Range('D49').Value = '''142.70'
$ws.Range('E49').Value = '  +1.04%  '
$ws.Range('D50').Value = '''0.0357'
$ws.Range('E50').Value = '  +1.16%  '
$ws.Range('D51').Value = '''25.77'
$ws.Range('E51').Value = '  +12.93%  '
